$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "gnk;kk"
$ws.Range("B2").Value = "bknkn"

$ws.Range("B2").Select()
